$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Inventario"

# Remove rows 5-7 (old RUT/ROL... extra rows), keep only 4 rows
$ws.Rows("5:7").Delete()

# Update remaining header/value cells for Inventario sheet
$ws.Range("A2").Value = "NOMBRE"
$ws.Range("B2").Value = "Patio trasero"
$ws.Range("A3").Value = "MÁXIMO STOCK"
$ws.Range("B3").Value = 200
$ws.Range("A4").Value = "ÚLTIMA ACTUALIZACIÓN"
$ws.Range("B4").Value = "28/10/2024"

# Add the new Productos sheet right after Inventario
$ws2 = $wb.Worksheets.Add($null, $ws)
$ws2.Name = "Productos"

# Copy the three header/row formats from Inventario onto the full 6-column width of Productos
$ws.Range("A1:B1").Copy()
$ws2.Range("A1:F1").PasteSpecial(-4122)
$ws.Range("A2:B2").Copy()
$ws2.Range("A2:F2").PasteSpecial(-4122)
$ws.Range("A3:B3").Copy()
$ws2.Range("A3:F3").PasteSpecial(-4122)
$ws.Range("A2:B2").Copy()
$ws2.Range("A4:F4").PasteSpecial(-4122)
$ws.Range("A3:B3").Copy()
$ws2.Range("A5:F5").PasteSpecial(-4122)
$ws.Range("A2:B2").Copy()
$ws2.Range("A6:F6").PasteSpecial(-4122)
$ws.Range("A3:B3").Copy()
$ws2.Range("A7:F7").PasteSpecial(-4122)
$ws.Range("A2:B2").Copy()
$ws2.Range("A8:F8").PasteSpecial(-4122)
$ws.Range("A3:B3").Copy()
$ws2.Range("A9:F9").PasteSpecial(-4122)
$ws.Range("A2:B2").Copy()
$ws2.Range("A10:F10").PasteSpecial(-4122)

# Match column widths (30) across all 6 columns of Productos
$ws2.Columns.Item(1).ColumnWidth = 29.17
$ws2.Columns.Item(2).ColumnWidth = 29.17
$ws2.Columns.Item(3).ColumnWidth = 29.17
$ws2.Columns.Item(4).ColumnWidth = 29.17
$ws2.Columns.Item(5).ColumnWidth = 29.17
$ws2.Columns.Item(6).ColumnWidth = 29.17

# Populate the product table values
$ws2.Range("A1").Value = "NOMBRE DEL PRODUCTO"
$ws2.Range("B1").Value = "MARCA"
$ws2.Range("C1").Value = "DESCRIPCIÓN"
$ws2.Range("D1").Value = "CATEGORÍA"
$ws2.Range("E1").Value = "TIPO"
$ws2.Range("F1").Value = "CANTIDAD"
$ws2.Range("A2").Value = "Cerveza Artesanal Actualizada test"
$ws2.Range("B2").Value = "NO REGISTRADO"
$ws2.Range("C2").Value = "holabrofdssdtest"
$ws2.Range("D2").Value = "Cerveza"
$ws2.Range("E2").Value = "Sin Alcohol"
$ws2.Range("F2").Value = 50
$ws2.Range("A3").Value = "holabro2dfsad3"
$ws2.Range("B3").Value = "NO REGISTRADO"
$ws2.Range("C3").Value = "testttasdasd"
$ws2.Range("D3").Value = "Vino"
$ws2.Range("E3").Value = "Sin Alcohol"
$ws2.Range("F3").Value = 20
$ws2.Range("A4").Value = "Ron 2 test"
$ws2.Range("B4").Value = "NO REGISTRADO"
$ws2.Range("C4").Value = "asasfafasfsafsfasasffasasfda"
$ws2.Range("D4").Value = "Ron"
$ws2.Range("E4").Value = "Alcohólico"
$ws2.Range("F4").Value = 5
$ws2.Range("A5").Value = "Vodka"
$ws2.Range("B5").Value = "NO REGISTRADO"
$ws2.Range("C5").Value = "Nuevo Vodka Vegano"
$ws2.Range("D5").Value = "Otro"
$ws2.Range("E5").Value = "Alcohólico"
$ws2.Range("F5").Value = 22
$ws2.Range("A6").Value = "nuevo test 2"
$ws2.Range("B6").Value = "NO REGISTRADO"
$ws2.Range("C6").Value = "testttasdasd"
$ws2.Range("D6").Value = "Agua mineral"
$ws2.Range("E6").Value = "Sin Alcohol"
$ws2.Range("F6").Value = 20
$ws2.Range("A7").Value = "NuevoTest25"
$ws2.Range("B7").Value = "NO REGISTRADO"
$ws2.Range("C7").Value = "holabrofdssdtest"
$ws2.Range("D7").Value = "Cigarrillo"
$ws2.Range("E7").Value = "Otro"
$ws2.Range("F7").Value = 25
$ws2.Range("A8").Value = "Cerveza Artesanal Actualizada tesasdasdt"
$ws2.Range("B8").Value = "NO REGISTRADO"
$ws2.Range("C8").Value = "testttasdasd"
$ws2.Range("D8").Value = "Snack"
$ws2.Range("E8").Value = "Otro"
$ws2.Range("F8").Value = 23
$ws2.Range("A9").Value = "Cerveza Artesanal Actualizada testasdasd"
$ws2.Range("B9").Value = "NO REGISTRADO"
$ws2.Range("C9").Value = "holabrofdssdtest"
$ws2.Range("D9").Value = "Néctar"
$ws2.Range("E9").Value = "Otro"
$ws2.Range("F9").Value = 25
$ws2.Range("A10").Value = "Ron"
$ws2.Range("B10").Value = "testqwrwqd"
$ws2.Range("C10").Value = "Ron Vodka cualificado por admins"
$ws2.Range("D10").Value = "Ron"
$ws2.Range("E10").Value = "Alcohólico"
$ws2.Range("F10").Value = 3
